$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header row columns: "<name>_old" -> "<name>_FV2304" (A:J) and
# "<name>_new" -> "<name>_FV2310" (L:U). Column K ("diff") stays unchanged.
$oldCols   = @("A","B","C","D","E","F","G","H","I","J")
$newCols   = @("L","M","N","O","P","Q","R","S","T","U")
$baseNames = @("Segmentname","Segmentgruppe","Segment","Datenelement","Segment ID","Code","Qualifier","Beschreibung","Bedingungsausdruck","Bedingung")

for ($i = 0; $i -lt $oldCols.Length; $i++) {
    $ws.Range($oldCols[$i] + "1").Value = "$($baseNames[$i])_FV2304"
}
for ($i = 0; $i -lt $newCols.Length; $i++) {
    $ws.Range($newCols[$i] + "1").Value = "$($baseNames[$i])_FV2310"
}

# Turn the whole used range into an Excel Table (Table1) with headers.
$range = $ws.Range("A1:U55")
$lo = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $range, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$lo.Name = "Table1"

# Freeze the header row (split below row 1, top-left of scrollable area is A2).
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
